$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definition")
$ws.Name = "Instrument Properties"
